# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to the newly scraped totals.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new F value
$sheetExhibit.Range("F3").Value  = 3798
$sheetExhibit.Range("F4").Value  = 2283
$sheetExhibit.Range("F5").Value  = 448
$sheetExhibit.Range("F7").Value  = 19
$sheetExhibit.Range("F8").Value  = 183
$sheetExhibit.Range("F10").Value = 91
$sheetExhibit.Range("F11").Value = 1413
$sheetExhibit.Range("F13").Value = 2377
$sheetExhibit.Range("F14").Value = 167

# Sheet "全部类型": same events, but at different row offsets because
# this sheet interleaves extra "演出" rows
$sheetAll.Range("F3").Value  = 3798
$sheetAll.Range("F4").Value  = 2283
$sheetAll.Range("F5").Value  = 448
$sheetAll.Range("F7").Value  = 19
$sheetAll.Range("F9").Value  = 183
$sheetAll.Range("F11").Value = 91
$sheetAll.Range("F14").Value = 1413
$sheetAll.Range("F16").Value = 2377
$sheetAll.Range("F17").Value = 167
